$wb = $excel.ActiveWorkbook

# --- Rename "Requested quantity" headers on the existing sheets ---
$wsWeekly = $wb.Worksheets.Item("Weekly Quantity")
$wsWeekly.Range("B1").Value = "Weekly_PO_Qty"

$wsMonthly = $wb.Worksheets.Item("Monthly Trend")
$wsMonthly.Range("B1").Value = "Monthly_PO_Qty"

# --- Add the new "PO Forecast" sheet after "Monthly Trend" ---
$wsForecast = $wb.Worksheets.Add($null, $wsMonthly)
$wsForecast.Name = "PO Forecast"

# Match the outline/page-setup and margin conventions used by the other sheets
$wsForecast.Outline.SummaryRow = 1
$wsForecast.Outline.SummaryColumn = 1
$wsForecast.PageSetup.LeftMargin = 54
$wsForecast.PageSetup.RightMargin = 54
$wsForecast.PageSetup.TopMargin = 72
$wsForecast.PageSetup.BottomMargin = 72
$wsForecast.PageSetup.HeaderMargin = 36
$wsForecast.PageSetup.FooterMargin = 36

# Headers
$wsForecast.Range("A1").Value = "ds"
$wsForecast.Range("B1").Value = "PO_Forecast"
$wsForecast.Range("C1").Value = "yhat_lower"
$wsForecast.Range("D1").Value = "yhat_upper"

# Data rows
$data = @(
    @(45578.99999999999, 3,   -21.19709484387107, 28.97635921882362),
    @(45592.99999999999, 69,   42.09846947689152, 96.06154746459546),
    @(45613.99999999999, 168, 143.4198888387816,  194.057014198614),
    @(45620.99999999999, 201, 175.4708509162121,  227.8299883965523),
    @(45627.99999999999, 234, 207.8616247595014,  258.4124669169067),
    @(45634.99999999999, 268, 241.1340417641401,  293.9644200429433),
    @(45641.99999999999, 301, 274.7920322287033,  327.6358518007075),
    @(45648.99999999999, 334, 307.6713289810671,  357.7216046712565),
    @(45655.99999999999, 367, 342.505516809932,   396.1084815132182),
    @(45662.99999999999, 400, 373.848430746858,   426.4630563329343),
    @(45669.99999999999, 433, 405.9395738938655,  459.6564283362459)
)

$r = 2
foreach ($row in $data) {
    $wsForecast.Cells.Item($r, 1).Value = $row[0]
    $wsForecast.Cells.Item($r, 2).Value = $row[1]
    $wsForecast.Cells.Item($r, 3).Value = $row[2]
    $wsForecast.Cells.Item($r, 4).Value = $row[3]
    $r++
}

# --- Formatting: match the header / date-column styles used on the other sheets ---
$wsWeekly.Range("A1:B1").Copy()
$wsForecast.Range("A1:D1").PasteSpecial(-4122)  # xlPasteFormats

$wsWeekly.Range("A2").Copy()
$wsForecast.Range("A2:A12").PasteSpecial(-4122)  # xlPasteFormats

$excel.CutCopyMode = 0

# Restore original active sheet/tab selection
$wsWeekly.Activate()
